$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without triggering Excel's
# "smart entry" auto-conversion (e.g. "85%" -> number 0.85 formatted as percent).
# We stage the literal text as a formula result in a scratch cell far outside the
# used range, copy it, and paste-special values-only into the target - this keeps
# the target cell's existing style/number-format untouched and leaves a plain text value.
function Set-LiteralText($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace("`"", "`"`"")
    $scratch.Formula = "=`"" + $escaped + "`""
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.ClearContents()
    $excel.CutCopyMode = $false
}

$ws.Range("E2").Value = "2026-02-11 07:48:35"
Set-LiteralText $ws.Range("H2") "85%"
$ws.Range("O2").Value = "1.9 °C"
$ws.Range("E3").Value = "2026-02-11 07:48:37"
Set-LiteralText $ws.Range("H3") "91%"
$ws.Range("E4").Value = "2026-02-11 07:48:39"
Set-LiteralText $ws.Range("H4") "72%"
$ws.Range("K4").Value = "0.0 MJ/m2"
$ws.Range("L4").Value = "67.3 km/h - 279º 7:29 TU"
$ws.Range("O4").Value = "13.0 °C"
$ws.Range("E5").Value = "2026-02-11 07:48:42"
Set-LiteralText $ws.Range("H5") "91%"
$ws.Range("E6").Value = "2026-02-11 07:48:45"
Set-LiteralText $ws.Range("H6") "94%"
$ws.Range("J6").Value = "1003.3 hPa"
$ws.Range("M6").Value = "13.0 °C 7:03 TU"
$ws.Range("O6").Value = "10.8 °C"
$ws.Range("E7").Value = "2026-02-11 07:48:47"
$ws.Range("J7").Value = "1002.7 hPa"
$ws.Range("K7").Value = "0.0 MJ/m2"
$ws.Range("N7").Value = "17.8 °C 7:25 TU"
$ws.Range("E8").Value = "2026-02-11 07:48:49"
$ws.Range("J8").Value = "1002.3 hPa"
$ws.Range("N8").Value = "13.7 °C 7:14 TU"
$ws.Range("O8").Value = "15.0 °C"
$ws.Range("E9").Value = "2026-02-11 07:48:51"
$ws.Range("K9").Value = "0.1 MJ/m2"
$ws.Range("O9").Value = "9.3 °C"
$ws.Range("E10").Value = "2026-02-11 07:48:54"
Set-LiteralText $ws.Range("H10") "91%"
$ws.Range("K10").Value = "0.0 MJ/m2"
$ws.Range("E11").Value = "2026-02-11 07:48:56"
$ws.Range("N11").Value = "2.5 °C 7:05 TU"
$ws.Range("O11").Value = "4.4 °C"
$ws.Range("E12").Value = "2026-02-11 07:48:59"
$ws.Range("O12").Value = "9.1 °C"
$ws.Range("E13").Value = "2026-02-11 07:49:01"
$ws.Range("J13").Value = "1005.6 hPa"
$ws.Range("K13").Value = "0.0 MJ/m2"
$ws.Range("E14").Value = "2026-02-11 07:49:04"
$ws.Range("N14").Value = "16.6 °C 7:04 TU"
$ws.Range("E15").Value = "2026-02-11 07:49:06"
$ws.Range("O15").Value = "9.2 °C"
$ws.Range("E16").Value = "2026-02-11 07:49:08"
Set-LiteralText $ws.Range("H16") "55%"
$ws.Range("E17").Value = "2026-02-11 07:49:10"
Set-LiteralText $ws.Range("H17") "88%"
$ws.Range("E18").Value = "2026-02-11 07:49:13"
Set-LiteralText $ws.Range("H18") "95%"
$ws.Range("M18").Value = "12.3 °C 7:29 TU"
$ws.Range("O18").Value = "9.4 °C"
$ws.Range("E19").Value = "2026-02-11 07:49:16"
$ws.Range("O19").Value = "8.1 °C"
$ws.Range("E20").Value = "2026-02-11 07:49:18"
$ws.Range("N20").Value = "-2.6 °C 7:03 TU"
$ws.Range("O20").Value = "-0.9 °C"
$ws.Range("E21").Value = "2026-02-11 07:49:20"
$ws.Range("J21").Value = "1006.2 hPa"
$ws.Range("N21").Value = "3.8 °C 7:08 TU"
$ws.Range("O21").Value = "5.6 °C"
$ws.Range("E22").Value = "2026-02-11 07:49:23"
$ws.Range("O22").Value = "-2.8 °C"
$ws.Range("E23").Value = "2026-02-11 07:49:25"
Set-LiteralText $ws.Range("H23") "66%"
$ws.Range("K23").Value = "0.0 MJ/m2"
$ws.Range("E24").Value = "2026-02-11 07:49:28"
$ws.Range("J24").Value = "1007.3 hPa"
$ws.Range("E25").Value = "2026-02-11 07:49:30"
Set-LiteralText $ws.Range("H25") "60%"
$ws.Range("K25").Value = "0.0 MJ/m2"
$ws.Range("E26").Value = "2026-02-11 07:49:33"
Set-LiteralText $ws.Range("H26") "83%"
$ws.Range("O26").Value = "4.1 °C"
$ws.Range("E27").Value = "2026-02-11 07:49:35"
Set-LiteralText $ws.Range("H27") "86%"
$ws.Range("E28").Value = "2026-02-11 07:49:38"
$ws.Range("J28").Value = "1004.2 hPa"
$ws.Range("O28").Value = "7.5 °C"
$ws.Range("E29").Value = "2026-02-11 07:49:40"
$ws.Range("N29").Value = "7.4 °C 7:26 TU"
$ws.Range("O29").Value = "9.7 °C"
$ws.Range("E30").Value = "2026-02-11 07:49:43"
$ws.Range("J30").Value = "1003.3 hPa"
$ws.Range("N30").Value = "7.8 °C 7:04 TU"
$ws.Range("O30").Value = "9.4 °C"
$ws.Range("E31").Value = "2026-02-11 07:49:45"
Set-LiteralText $ws.Range("H31") "68%"
$ws.Range("J31").Value = "1002.3 hPa"
$ws.Range("K31").Value = "0.0 MJ/m2"
$ws.Range("E32").Value = "2026-02-11 07:49:48"
$ws.Range("K32").Value = "0.0 MJ/m2"
$ws.Range("L32").Value = "43.2 km/h - 312º 7:24 TU"
$ws.Range("E33").Value = "2026-02-11 07:49:50"
$ws.Range("J33").Value = "1005.5 hPa"
$ws.Range("N33").Value = "2.4 °C 7:01 TU"
$ws.Range("O33").Value = "4.4 °C"
$ws.Range("E34").Value = "2026-02-11 07:49:53"
Set-LiteralText $ws.Range("H34") "71%"
$ws.Range("E35").Value = "2026-02-11 07:49:55"
Set-LiteralText $ws.Range("H35") "63%"
$ws.Range("J35").Value = "1008.6 hPa"
$ws.Range("E36").Value = "2026-02-11 07:49:58"
$ws.Range("J36").Value = "1003.4 hPa"
$ws.Range("N36").Value = "8.1 °C 7:21 TU"
$ws.Range("O36").Value = "10.3 °C"
$ws.Range("E37").Value = "2026-02-11 07:50:00"
$ws.Range("J37").Value = "1005.5 hPa"
$ws.Range("L37").Value = "23.4 km/h - 243º 7:14 TU"
$ws.Range("N37").Value = "5.2 °C 7:20 TU"
$ws.Range("O37").Value = "7.0 °C"
$ws.Range("E38").Value = "2026-02-11 07:50:03"
Set-LiteralText $ws.Range("H38") "76%"
$ws.Range("K38").Value = "0.0 MJ/m2"
$ws.Range("O38").Value = "12.7 °C"
$ws.Range("E39").Value = "2026-02-11 07:50:05"
Set-LiteralText $ws.Range("H39") "61%"
$ws.Range("E40").Value = "2026-02-11 07:50:08"
$ws.Range("J40").Value = "1007.6 hPa"
$ws.Range("O40").Value = "4.6 °C"
$ws.Range("E41").Value = "2026-02-11 07:50:10"
$ws.Range("J41").Value = "1003.7 hPa"
$ws.Range("K41").Value = "0.0 MJ/m2"
$ws.Range("N41").Value = "17.5 °C 7:28 TU"
$ws.Range("O41").Value = "19.1 °C"
$ws.Range("E42").Value = "2026-02-11 07:50:13"
$ws.Range("N42").Value = "7.2 °C 7:05 TU"
$ws.Range("O42").Value = "9.7 °C"
$ws.Range("E43").Value = "2026-02-11 07:50:16"
$ws.Range("N43").Value = "10.5 °C 7:00 TU"
$ws.Range("O43").Value = "12.3 °C"
$ws.Range("E44").Value = "2026-02-11 07:50:18"
Set-LiteralText $ws.Range("H44") "72%"
$ws.Range("E45").Value = "2026-02-11 07:50:21"
$ws.Range("J45").Value = "1008.0 hPa"
$ws.Range("L45").Value = "18.4 km/h - 101º 7:24 TU"
$ws.Range("E46").Value = "2026-02-11 07:50:23"
$ws.Range("J46").Value = "1007.8 hPa"
$ws.Range("N46").Value = "15.7 °C 7:29 TU"
$ws.Range("O46").Value = "17.8 °C"
